$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the existing "sum" header (G1) onto the new "Save" header (H1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
